$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '24.508.40'
$ws.Range('E2').Value = '  -1.06%  '
$ws.Range('D3').Value = '1.692.72'
$ws.Range('E3').Value = '  -0.69%  '
$ws.Range('D4').Value = '1.001'
$ws.Range('E4').Value = '  -0.25%  '
$ws.Range('D5').Value = '316.12'
$ws.Range('E5').Value = '  -0.26%  '
$ws.Range('D6').Value = '1.001'
$ws.Range('E6').Value = '  -0.29%  '
$ws.Range('D7').Value = '0.3903'
$ws.Range('E7').Value = '  -0.98%  '
$ws.Range('D8').Value = '0.4046'
$ws.Range('E8').Value = '  -0.05%  '
$ws.Range('E9').Value = '  -2.12%  '
$ws.Range('D10').Value = '1.002'
$ws.Range('E10').Value = '  -0.06%  '
$ws.Range('D11').Value = '52.82'
$ws.Range('E11').Value = '  -1.25%  '
$ws.Range('D12').Value = '0.08772'
$ws.Range('E12').Value = '  -1.31%  '
$ws.Range('D13').Value = '26.71'
$ws.Range('E13').Value = '  +12.60%  '
$ws.Range('D14').Value = '7.474'
$ws.Range('E14').Value = '  +1.71%  '
$ws.Range('D15').Value = '8.142'
$ws.Range('E15').Value = '  +0.98%  '
$ws.Range('D16').Value = '0.00001347'
$ws.Range('E16').Value = '  +1.42%  '
$ws.Range('D17').Value = '1.685.99'
$ws.Range('E17').Value = '  -1.90%  '
$ws.Range('D18').Value = '98.02'
$ws.Range('E18').Value = '  -1.98%  '
$ws.Range('D19').Value = '0.07184'
$ws.Range('E19').Value = '  +1.92%  '
$ws.Range('D20').Value = '20.39'
$ws.Range('E20').Value = '  +3.22%  '
$ws.Range('D21').Value = '7.299'
$ws.Range('E21').Value = '  +3.10%  '
$ws.Range('D22').Value = '1.000'
$ws.Range('E22').Value = '  -0.09%  '
$ws.Range('D23').Value = '14.29'
$ws.Range('E23').Value = '  -1.34%  '
$ws.Range('D24').Value = '24.493.02'
$ws.Range('E24').Value = '  -1.09%  '
$ws.Range('D25').Value = '3.021'
$ws.Range('E25').Value = '  -6.22%  '
$ws.Range('E26').Value = '  -1.28%  '
$ws.Range('D27').Value = '22.63'
$ws.Range('E27').Value = '  -0.73%  '
$ws.Range('D28').Value = '167.38'
$ws.Range('E28').Value = '  +3.10%  '
$ws.Range('D29').Value = '8.418'
$ws.Range('E29').Value = '  -4.70%  '
$ws.Range('D30').Value = '5.396'
$ws.Range('E30').Value = '  +4.27%  '
$ws.Range('D31').Value = '138.24'
$ws.Range('E31').Value = '  +1.24%  '
$ws.Range('D32').Value = '1.872.66'
$ws.Range('E32').Value = '  -1.13%  '
$ws.Range('B33').Value = 'WEMIXTOKEN'
$ws.Range('C33').Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range('D33').Value = '2.157'
$ws.Range('E33').Value = '  +8.64%  '
$ws.Range('B34').Value = 'Hedera'
$ws.Range('C34').Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range('D34').Value = '0.08717'
$ws.Range('E34').Value = '  -2.32%  '
$ws.Range('B35').Value = 'InternetComputer(DFINITY)'
$ws.Range('C35').Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range('D35').Value = '7.277'
$ws.Range('E35').Value = '  -8.77%  '
$ws.Range('D36').Value = '1.037'
$ws.Range('E36').Value = '  -4.24%  '
$ws.Range('E37').Value = '  +6.66%  '
$ws.Range('D38').Value = '0.2783'
$ws.Range('E38').Value = '  +0.68%  '
$ws.Range('D39').Value = '10.89'
$ws.Range('E39').Value = '  -1.80%  '
$ws.Range('B40').Value = 'Stellar'
$ws.Range('C40').Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range('D40').Value = '0.09148'
$ws.Range('E40').Value = '  -0.46%  '
$ws.Range('B41').Value = 'TheSandbox'
$ws.Range('C41').Value = 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
$ws.Range('D41').Value = '0.8044'
$ws.Range('E41').Value = '  +4.04%  '
$ws.Range('D42').Value = '14.18'
$ws.Range('E42').Value = '  -3.18%  '
$ws.Range('D43').Value = '1.477'
$ws.Range('E43').Value = '  +1.04%  '
$ws.Range('D44').Value = '17.45'
$ws.Range('E44').Value = '  +9.15%  '
$ws.Range('D45').Value = '2.666'
$ws.Range('E45').Value = '  +3.19%  '
$ws.Range('D46').Value = '0.7247'
$ws.Range('E46').Value = '  +0.37%  '
$ws.Range('D47').Value = '4.259'
$ws.Range('E47').Value = '  +1.09%  '
$ws.Range('D48').Value = '1.404'
$ws.Range('E48').Value = '  +4.95%  '
$ws.Range('D49').Value = '1.000'
$ws.Range('E49').Value = '  -0.17%  '
$ws.Range('D50').Value = '139.67'
$ws.Range('E50').Value = '  -0.96%  '
$ws.Range('D51').Value = '0.08162'
$ws.Range('E51').Value = '  +2.08%  '
